$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$v1 = @'
['HMV23/000055/0923/2', 'HMV23/000059/0923/5', 'HMV23/000055/0923/6', 'HMV23/000043/0823/12']
'@
$ws.Range("B4").Value = $v1

$v2 = @'
['DURING INSPECTION OBSERVED LH WING SLAT#5 WEATHER SEAL ERODED.SAME TO BE REPLACED.', 'DURING ARRIVAL INSPECTION , OBSERVED FOLLOWINS DICREPANCIES.\n\n1) RH SIDE WING NO#1 SLAT, NO.3 TRACK SQURE SEAL UPPER BULB SEAL DEGRADED.\n2) RH SIIDE WING SLAT NO#5 OUT BOARD WEATHER SEAL ERODED. SAME TO BE REPLACED.', 'DURING INSPECTION, OBSERVED RH SIDE WING SLAT NO#5 WEATHER SEAL ERODED.SAME TO BE REPLACED.', 'WHILE ARRIVAL INSPECTION OBSERVED RH WING SLAT#5 OUTBOARD WEATHER SEAL ERODED.SAME TO BE REPLACED.']
'@
$ws.Range("C4").Value = $v2

$v3 = @'
[2.0, 0.0, 3.0, 2.0]
'@
$ws.Range("E4").Value = $v3

$v4 = @'
['HMV23/000059/0923/6', 'HMV23/000055/0923/11']
'@
$ws.Range("B5").Value = $v4

$v5 = @'
['DURING INSPECTION FOUND FOLLOWING FORWARD CARGO FLOOR PANEL DAMAGED , TO BE FABRICATED/REPLACED \n\n\n\n131EF P/N D5367521000000 \n\n131FF P/N D5367400800000 \n\n131BF P/N D5367401200000 \n\n131JF P/N D5367400600000', 'DURING INSPECTION FOUND FOLLOWING FORWARD CARGO FLOOR PANEL DAMAGED , TO BE FABRICATED/REPLACED \n\n131 AF P/N D5367520400000\n\n131EF P/N D5367521000000\n\n131FF P/N D5367400800000\n\n131BF P/N D5367401200000\n\n131CF P/N D5367401100000\n\n131GF P/N D5367400700000']
'@
$ws.Range("C5").Value = $v5

$v6 = @'
['HMV23/000059/0923/8', 'HMV23/000055/0923/13', 'HMV23/000059/0923/9', 'HMV23/000055/0923/14']
'@
$ws.Range("B7").Value = $v6

$v7 = @'
['DURING ARRIVAL INSPECTION OF FORWARD CARGO , FOUND HARDWARE DAMAGED/MISSING', 'DURING ARRIVAL INSPECTION OF FORWARD CARGO , FOUND \n1) CARGO SECTION HOSE MISSING/ DAMAGED \n2) FOUND HARDWARE DAMAGED/MISSING', 'DURING ARRIVAL INSPECTION OF AFT CARGO , FOUND \n1) CARGO SECTION HOSE MISSING/ DAMAGED \n2) FOUND HARDWARE DAMAGED/MISSING', 'DURING ARRIVAL INSPECTION OF AFT CARGO , FOUND \n1) CARGO SECTION HOSE MISSING/ DAMAGED \n2) FOUND HARDWARE DAMAGED/MISSING']
'@
$ws.Range("C7").Value = $v7

$v8 = @'
[5.0, 5.0, 7.0, 5.0]
'@
$ws.Range("E7").Value = $v8

$v9 = @'
['HMV23/000059/0923/12', 'HMV23/000055/0923/24']
'@
$ws.Range("B10").Value = $v9

$v10 = @'
['DURING INSPECTION FOUND COMPRESSOR INLET FILTER BLOCKED , TO BE\n REPLACED', 'DURING INSPECTION FOUND COMPRESSOR INLET FILTER BLOCKED , TO BE REPLACED\nP/N 7006-15']
'@
$ws.Range("C10").Value = $v10

$v11 = @'
[3.0, 2.5]
'@
$ws.Range("E10").Value = $v11

$v12 = @'
['AIRCRAFT ARRIVAL INSPECTION', 'AWR23/003086/0923']
'@
$ws.Range("A11").Value = $v12

$v13 = @'
['HMV23/000055/0923/10', 'HMV23/000059/0923/15']
'@
$ws.Range("B11").Value = $v13

$v14 = @'
['DURING WEEKLY INSPECTION OBSERVED FIN: 11WL EPSU BATTERY INOP.', 'DURING ARRIVAL CHECK FOUND EPSU BATTERIES QTY#2 FIN:10WL AND FIN: 11WL INOP.']
'@
$ws.Range("C11").Value = $v14

$v15 = @'
[5.0, 6.0]
'@
$ws.Range("E11").Value = $v15

$v16 = @'
['HMV23/000055/0923/8', 'HMV23/000055/0923/1', 'HMV23/000055/0923/7', 'HMV23/000055/0923/9']
'@
$ws.Range("B14").Value = $v16

$v17 = @'
['DURING ARRIVAL INSPECTION, OBSERVED SMALL DENTS ON SLAT NO #1 TOP SKIN AT RH SIDE WING. LOCATION OF DENT 3.5 INCH FROM T/E AND 37 INCH FROM INBD EDGE OF SLAT.\nSLAT NO#1 P/N- D5746091000300\n S/N- SA7205917', 'DURING ARRIVAL INSPECTION , OBSERVED DENT ON SLAT NO #2 LIP AREA AT RH SIDE WING. ASSESSMENT TO BE CARRY OUT.\n\nSLAT NO #2 P/N- D5746092000500\n SA7205698', 'DURING ARRIVAL INSPECTION , OBSERVED DENT ON SLAT NO #2 LIP AREA AT RH SIDE WING. DENT LOCATION, 3 INCH FROM T/E AND 67.5 INCH FROM OUTBD EDGE .ASSESSMENT TO BE CARRY OUT.\nSLAT NO #2 P/N- D5746092000500\n S/N- SA7205698', 'DURING ARRIVAL INSPECTION, OBSERVED SMALL DENTS ON SLAT NO #1 TOP SKIN AT RH SIDE WING. LOCATION OF DENT 3.5 INCH FROM T/E AND 37 INCH FROM INBD EDGE OF SLAT.\nSLAT NO#1 P/N- D5746091000300\n S/N- SA7205917']
'@
$ws.Range("C14").Value = $v17

$v18 = @'
[3.0, 2.0, 7.0, 7.0]
'@
$ws.Range("E14").Value = $v18

$v19 = @'
['HMV23/000055/0923/5', 'HMV23/000055/0923/3', 'HMV23/000055/0923/4']
'@
$ws.Range("B15").Value = $v19

$v20 = @'
['DURING ARRIVAL INSPECTION FOUND FOLLOWING FAILURE MSG IN PFR:\n1. 38-31-41 TOILET ASSY LAV F\nFURTHER RECTIFICATION TO BE CARRIED OUT.', 'DURING ARRIVAL INSPECTION FOUND FOLLOWING WARNING MSG IN PFR:\n1. AUTO FLT AP OFF\nFURTHER RECTIFICATION TO BE CARRIED OUT.', 'DURING ARRIVAL INSPECTION FOUND FOLLOWING FAILURE MSG IN PFR:\n1. 34-53-31 ADF1(2RP1)\nFURTHER RECTIFICATION TO BE CARRIED OUT.']
'@
$ws.Range("C15").Value = $v20

$v21 = @'
[4.0, 2.0, 2.5]
'@
$ws.Range("E15").Value = $v21

$v22 = @'
['HMV23/000055/0923/19', 'HMV23/000055/0923/20']
'@
$ws.Range("B16").Value = $v22

$v23 = @'
['1.CARRY OUT HFEC INSPECTION AS PER NTM 51-10-08 ON DENT AND SURROUNDING AREA OF DENT ON RH WING SLAT#2 NOSE TOP SKIN TO ENSURE CRACK FREE CONDITION EXISTS AS REFERRED BY SRM TASK 57-40-00-283-007 REV:37 DTD:AUG 01/2023.\n\nDAMAGE DIMENSION AND LOCATION ARE\n\nDENT ON NOSE TOP SKIN OF SLAT #2 . DAMAGE DIMENSION \nA) LENGTH: 84MM\n\nB) WIDTH: 35MM\n\nC) DEPTH: 1.35MM\n\nD) LOCATION: DISTANCE FROM RIB WFX5997 FASTENER ROW CENTER =76MM, DISTANCE FROM UPPER GIRDER RIVET ROW CENTER =65MM', '1.CARRY OUT HFEC INSPECTION AS PER NTM 51-10-08 ON DENT AND SURROUNDING AREA OF DENT ON RH WING SLAT#1 NOSE TOP SKIN TO ENSURE CRACK FREE CONDITION EXISTS AS REFERRED BY SRM TASK 57-40-00-283-007 REV:37 DTD:AUG 01/2023.\n\nDAMAGE DIMENSION AND LOCATION ARE\n\nDENT ON NOSE TOP SKIN OF SLAT #1 . DAMAGE DIMENSION \nA) LENGTH: 28MM\n\nB) WIDTH: 25MM\n\nC) DEPTH: 0.30MM\n\nD) LOCATION: DISTANCE FROM RIB TRACK 2 RIVET ROW CENTER =109MM, DISTANCE FROM UPEER GIRDER RIVET ROW CENTER =25MM']
'@
$ws.Range("C16").Value = $v23

$v24 = @'
['HMV23/000043/0823/9', 'HMV23/000043/0823/8']
'@
$ws.Range("B19").Value = $v24

$v25 = @'
['POST REPAIR INSTALLATION OF LH WING SLAT #2(P/N-D5746092000400, S/N-SA7202688) TO BE CARRIED OUT .', 'REMOVAL OF LH WING SLAT #2 (P/N-D5746092000400, S/N-SA7202688)TO BE CARRIED OUT FOR REPAIR.']
'@
$ws.Range("C19").Value = $v25

$v26 = @'
[10.0, 5.0]
'@
$ws.Range("E19").Value = $v26

